$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B5").Value = "SingleUseId4"
$ws.Range("C5").Value = "Small"
$ws.Range("D5").Value = "Left"
$ws.Range("E5").Value = "LTR"
$ws.Range("F5").Value = "<>"
